# Daily refresh of the cryptos price list (GitHub Actions scheduled job).
# Updates Price (D) / Volume(1h) (E) text columns, and for row 47/48 also
# swaps which coin (Injective Protocol / Sui Network) occupies which row.
# All D/E values are stored as literal text (matching the sheet's existing
# inline-string cells) - a leading "'" is used where the new value would
# otherwise be auto-parsed as a number by Excel, which would silently
# re-format it (e.g. "32.90" -> "32.9") or lose precision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.783.26"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "3.702.45"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'678.07"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").Value = "'162.38"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").Value = "'0.148"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("D10").Value = "'7.12"
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("D13").Value = "'32.90"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "3.703.75"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "69.783.43"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("E17").Value = "  +2.09%  "
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").Value = "'472.79"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").Value = "'80.54"
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("D23").Value = "3.852.09"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("D24").Value = "'0.0000128"
$ws.Range("E24").Value = "  +3.15%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'11.05"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("D27").Value = "'9.15"
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("E30").Value = "  +1.76%  "
$ws.Range("D31").Value = "'6.63"
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").Value = "'26.90"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  +3.44%  "
$ws.Range("D35").Value = "3.691.51"
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("D36").Value = "'8.57"
$ws.Range("E36").Value = "  +4.69%  "
$ws.Range("D37").Value = "'6.19"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("D42").Value = "'169.05"
$ws.Range("E42").Value = "  +2.19%  "
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("E45").Value = "  +2.30%  "
$ws.Range("D46").Value = "'0.000282"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("B47").Value = "SuiNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D47").Value = "'1.11"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'28.03"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "'7.93"
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("E51").Value = "  +2.72%  "
